$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: '8Views'
$ws.Range("A2").Value = '8Views'
$ws.Range("B2").Value = 'http://www.8views.com/'
$ws.Range("C2").Value = '8Views is a fast-paced, end-to-end digital marketing company providing solutions across digital platforms. The core focus of the company is to help brands meet their business goals through the rapidly growing online space. Our solutions include search engine optimization, social media marketing, email marketing, ad campaigns, content marketing, analytics, and more.'
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 84
$ws.Range("F2").Value = 'Hiring since March 2023'
$ws.Range("G2").Value = $null
$ws.Range("H2").Value = 'Digital Advertising, Digital Marketing, Facebook Ads, Google AdWords, Google Analytics, Instagram Marketing, Search Engine Marketing (SEM), Certificate,  Letter of recommendation,  Informal dress code,  5 days a week,  Free snacks & beverages'
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 'Secunderabad, Hyderabad, Madhapur, Telangana'
$ws.Range("K2").Value = 'Stipend'
$ws.Range("L2").Value = '₹ 8,000 - 15,000 /month'
$ws.Range("M2").Value = $null
$ws.Range("N2").Value = 'https://internshala.com/internship/detail/performance-marketing-paid-ads-internship-in-multiple-locations-at-8views1744015239'

# Row 3: 'MentorBoxx'
$ws.Range("A3").Value = 'MentorBoxx'
$ws.Range("B3").Value = $null
$ws.Range("C3").Value = 'Our sole aim at MentorBoxx is to bridge the gap between universities & industries. We select 30 students every month to regularly interact with the right industry experts, work on live industry projects, and grasp as much industry knowledge as possible.'
$ws.Range("D3").Value = 13
$ws.Range("E3").Value = 1630
$ws.Range("F3").Value = 'Hiring since April 2021'
$ws.Range("G3").Value = $null
$ws.Range("H3").Value = 'Creative Writing, Email Marketing, English Proficiency (Spoken), English Proficiency (Written), Facebook Marketing, Google AdWords, Google Analytics, Instagram Marketing, Search Engine Marketing (SEM), Search Engine Optimization (SEO), Certificate,  Flexible work hours'
$ws.Range("I3").Value = 208
$ws.Range("J3").Value = 'Chennai, Coimbatore, Delhi, Gurgaon, Lucknow, Patna, Pune, Ranchi, Hyderabad, Mumbai, Varanasi, Jaipur, Noida, Bangalore, Andhra Tharhi'
$ws.Range("K3").Value = 'Stipend'
$ws.Range("L3").Value = '₹ 10,000 /month'
$ws.Range("M3").Value = $null
$ws.Range("N3").Value = 'https://internshala.com/internship/detail/part-time-digital-marketing-internship-in-multiple-locations-at-mentorboxx1744029239'

# Row 4: 'Vitals7'
$ws.Range("A4").Value = 'Vitals7'
$ws.Range("B4").Value = 'http://vitals7.com'
$ws.Range("C4").Value = 'Vitals7 is an innovative telehealth and AI-driven self-health monitoring platform. We empower users with cutting-edge technology, predictive analytics, and holistic healthcare solutions. Our platform integrates IoT-based health monitoring, AI-powered analytics, and digital health consultations to make healthcare more accessible and affordable. We are looking for passionate and creative digital marketing interns to join our team and help execute AI-powered marketing campaigns using the latest automation tools.'
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 'Hiring since March 2025'
$ws.Range("G4").Value = $null
$ws.Range("H4").Value = 'Django, Flask, Machine Learning, Python, Certificate,  Letter of recommendation,  Flexible work hours'
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 'Hyderabad, Bangalore                                                    (Hybrid)'
$ws.Range("K4").Value = 'Stipend'
$ws.Range("L4").Value = '₹ 10,000 - 40,000 /month'
$ws.Range("M4").Value = $null
$ws.Range("N4").Value = 'https://internshala.com/internship/detail/python-development-internship-in-multiple-locations-at-vitals71743684810'

# Row 5: 'Tex N Co'
$ws.Range("A5").Value = 'Tex N Co'
$ws.Range("B5").Value = $null
$ws.Range("C5").Value = 'We are an interior design and construction firm handling design and execution for both residential and commercial interior, construction, and renovation projects.'
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 'Hiring since July 2021'
$ws.Range("G5").Value = $null
$ws.Range("H5").Value = 'AutoCAD, Data Analytics, Engineering Drawing, Engineering Surveying, English Proficiency (Written), Google Docs, Google Sheets, Google SketchUp , Certificate,  Letter of recommendation,  Informal dress code,  Free snacks & beverages'
$ws.Range("I5").Value = 106
$ws.Range("J5").Value = 'Hyderabad'
$ws.Range("K5").Value = 'Stipend'
$ws.Range("L5").Value = '₹ 10,000 /month'
$ws.Range("M5").Value = $null
$ws.Range("N5").Value = 'https://internshala.com/internship/detail/technical-assistant-internship-in-hyderabad-at-tex-n-co1743582036'

# Row 6: 'J K Arts'
$ws.Range("A6").Value = 'J K Arts'
$ws.Range("B6").Value = $null
$ws.Range("C6").Value = 'We are a startup, mainly focused on sentimental analysis, PR activities, and reputation management.'
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 'Hiring since March 2025'
$ws.Range("G6").Value = $null
$ws.Range("H6").Value = 'Adobe After Effects, Adobe Illustrator, Adobe Photoshop, Adobe Premiere Pro, Business Development, Client Relationship, Conflict Management, Content Management, Critical thinking, Data Analysis, Data Extraction, Data Manipulation, Final Cut Pro, Influencer Marketing, Problem Solving, Report Generation, Research and Analytics, Resource Management, Video Editing, Videography, Certificate'
$ws.Range("I6").Value = 32
$ws.Range("J6").Value = 'Hyderabad'
$ws.Range("K6").Value = 'Stipend'
$ws.Range("L6").Value = '₹ 8,000 - 12,000 /month'
$ws.Range("M6").Value = $null
$ws.Range("N6").Value = 'https://internshala.com/internship/detail/associate-consultant-internship-in-hyderabad-at-j-k-arts1743423990'

# Row 7: 'The Affordable Organic Store'
$ws.Range("A7").Value = 'The Affordable Organic Store'
$ws.Range("B7").Value = 'https://theaffordableorganicstore.com/'
$ws.Range("C7").Value = 'We are a bunch of enthusiasts who want to make good food affordable again by cutting out the middlemen. We want to set an example that a sustainable business can also be run without exploiting the consumers with high prices. We will be successful if everyone in the country starts doing what we do. We will change the rules of retail and modern trade. We will bring the power back to the producers from the retailers or middlemen. From our team, we want only the best. We are a tech-driven company with a focus on constant innovation. We need smart people with the best analytical and communication skills and a great heart. If we succeed, we will create a world order where everyone has access to good food and will lead a content life. We promise that we will take care of you at all times.'
$ws.Range("D7").Value = 156
$ws.Range("E7").Value = 714
$ws.Range("F7").Value = 'Hiring since May 2020'
$ws.Range("G7").Value = $null
$ws.Range("H7").Value = 'Data Analytics, Data Science, MS-Excel, MS-Office, Power BI, SQL, Certificate,  Letter of recommendation'
$ws.Range("I7").Value = 226
$ws.Range("J7").Value = 'Hyderabad'
$ws.Range("K7").Value = 'Stipend'
$ws.Range("L7").Value = '₹ 5,000 /month'
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = 'https://internshala.com/internship/detail/data-analytics-internship-in-hyderabad-at-the-affordable-organic-store1742636460'

# Row 8: 'RIAI'
$ws.Range("A8").Value = 'RIAI'
$ws.Range("B8").Value = 'https://riai.co.in/'
$ws.Range("C8").Value = 'At RIAI, we are not just developing AI; we are redefining how AI interacts with users and businesses through continuous research and development.'
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 'Hiring since March 2025'
$ws.Range("G8").Value = $null
$ws.Range("H8").Value = 'Email Marketing, Facebook Ads, Google AdWords, Google Analytics, Search Engine Optimization (SEO), WordPress, Certificate,  Letter of recommendation'
$ws.Range("I8").Value = 65
$ws.Range("J8").Value = 'Hyderabad'
$ws.Range("K8").Value = 'Stipend'
$ws.Range("L8").Value = '₹ 18,000 /month'
$ws.Range("M8").Value = $null
$ws.Range("N8").Value = 'https://internshala.com/internship/detail/psychology-behavioral-marketing-strategic-communications-internship-in-hyderabad-at-riai1742462010'

# Row 9: 'Grow Easy Hair Care Solutions Private Limited'
$ws.Range("A9").Value = 'Grow Easy Hair Care Solutions Private Limited'
$ws.Range("B9").Value = $null
$ws.Range("C9").Value = 'We are a startup in the hair care industry working to overcome problems of traditional hair care methods and provide chemical-free products to customers.'
$ws.Range("D9").Value = 12
$ws.Range("E9").Value = 110
$ws.Range("F9").Value = 'Hiring since March 2021'
$ws.Range("G9").Value = $null
$ws.Range("H9").Value = 'Data Analytics, Research and Analytics, Certificate,  Letter of recommendation,  Informal dress code,  Free snacks & beverages,  Job offer'
$ws.Range("I9").Value = 86
$ws.Range("J9").Value = 'Hyderabad'
$ws.Range("K9").Value = 'Stipend'
$ws.Range("L9").Value = '₹ 2,000 /month'
$ws.Range("M9").Value = $null
$ws.Range("N9").Value = 'https://internshala.com/internship/detail/business-analytics-internship-in-hyderabad-at-grow-easy-hair-care-solutions-private-limited1742447585'

# Row 10: 'APTAGRIM CONSULTING PRIVATE LIMITED'
$ws.Range("A10").Value = 'APTAGRIM CONSULTING PRIVATE LIMITED'
$ws.Range("B10").Value = 'https://aptagrim.com/'
$ws.Range("C10").Value = 'APTAGRIM CONSULTING PRIVATE LIMITED is a DeepTech AI company that offers a range of AI services to businesses across different industries. With expertise in deep learning, CNN, computer vision, NLP, and chatbot development, Aptagrim provides AI-powered solutions that automate business processes, analyze data, and enhance customer engagement. The company also offers data engineering and business intelligence services that turn data into actionable insights. Aptagrim is equipped to handle product engineering services, incubate ideas, and develop MVPs for startups and entrepreneurs.'
$ws.Range("D10").Value = 21
$ws.Range("E10").Value = 124
$ws.Range("F10").Value = 'Hiring since March 2023'
$ws.Range("G10").Value = $null
$ws.Range("H10").Value = 'Computer Vision, Deep Learning, Natural Language Processing (NLP), Certificate,  Letter of recommendation,  5 days a week,  Free snacks & beverages'
$ws.Range("I10").Value = 674
$ws.Range("J10").Value = 'Hyderabad'
$ws.Range("K10").Value = 'Stipend'
$ws.Range("L10").Value = '₹ 15,000 /month'
$ws.Range("M10").Value = $null
$ws.Range("N10").Value = 'https://internshala.com/internship/detail/ai-engineer-internship-in-hyderabad-at-aptagrim-consulting-private-limited1742035373'

# Row 11: 'Denary Media Private Limited'
$ws.Range("A11").Value = 'Denary Media Private Limited'
$ws.Range("B11").Value = $null
$ws.Range("C11").Value = 'We are a relatively new company based out of Secunderabad. Our founders have been in the digital space for over 6 years and are quite adept at delivering data-driven results. We are currently working with two large hospital chains and a few other large brands in Hyderabad.'
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 30
$ws.Range("F11").Value = 'Hiring since October 2019'
$ws.Range("G11").Value = $null
$ws.Range("H11").Value = 'Digital Marketing, English Proficiency (Written), Google Analytics, Search Engine Optimization (SEO), Certificate'
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 'Hyderabad'
$ws.Range("K11").Value = 'Stipend'
$ws.Range("L11").Value = '₹ 10,000 /month'
$ws.Range("M11").Value = $null
$ws.Range("N11").Value = 'https://internshala.com/internship/detail/search-engine-optimization-seo-internship-in-hyderabad-at-denary-media-private-limited1741765638'

# Row 12: 'Megaminds IT Services'
$ws.Range("A12").Value = 'Megaminds IT Services'
$ws.Range("B12").Value = $null
$ws.Range("C12").Value = 'We provide application solutions and development needs to clients. We offer business solutions for computing and creative excellence and deliver innovative and cost-effective solutions with ethics. Our experienced and committed team, with its strong focus on technology, is the backbone of our company in delivering the synergy of creative solutions. We design, develop, and deliver cost-effective and high-quality software applications. We provide e-commerce, retail, manufacturing, and many other services. We work on political campaigning with transparency. Megaminds IT Services (now Megaminds IT & Job Consultancy Services) has started job consultancy services too.'
$ws.Range("D12").Value = 81
$ws.Range("E12").Value = 304
$ws.Range("F12").Value = 'Hiring since March 2019'
$ws.Range("G12").Value = $null
$ws.Range("H12").Value = 'Algorithms, Business Analysis, Computer Vision, Data Analysis, Data Analytics, Data Science, English Proficiency (Spoken), English Proficiency (Written), LaTeX, Machine Learning, MS-PowerPoint, MS-Word, Power BI, Research and Analytics, Certificate,  Letter of recommendation'
$ws.Range("I12").Value = 178
$ws.Range("J12").Value = 'Hyderabad'
$ws.Range("K12").Value = 'Stipend'
$ws.Range("L12").Value = '₹ 6,000-8,000 /month'
$ws.Range("M12").Value = $null
$ws.Range("N12").Value = 'https://internshala.com/internship/detail/content-writer-internship-in-hyderabad-at-megaminds-it-services1741669406'

# Row 13: 'LYWO Recruitment Consulting LLP'
$ws.Range("A13").Value = 'LYWO Recruitment Consulting LLP'
$ws.Range("B13").Value = $null
$ws.Range("C13").Value = 'LYWO is an early-stage start-up with the objective to identify or develop an AI-assisted behavioral model that is simple enough to be used by an organization of any size.'
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = 'Hiring since March 2025'
$ws.Range("G13").Value = $null
$ws.Range("H13").Value = 'Adobe After Effects, Adobe Creative Suite, Adobe Illustrator, Adobe InDesign, Adobe Photoshop, Adobe Photoshop Lightroom CC, Adobe Premiere Pro, Computer Vision, Visual Basic (VB), Certificate,  Letter of recommendation,  Job offer'
$ws.Range("I13").Value = 134
$ws.Range("J13").Value = 'Hyderabad'
$ws.Range("K13").Value = 'Stipend'
$ws.Range("L13").Value = '₹ 15,000 /month'
$ws.Range("M13").Value = $null
$ws.Range("N13").Value = 'https://internshala.com/internship/detail/graphic-design-internship-in-hyderabad-at-lywo-recruitment-consulting-llp1741585446'

# Row 14: 'The Short Media'
$ws.Range("A14").Value = 'The Short Media'
$ws.Range("B14").Value = 'https://theshortmedia.com/'
$ws.Range("C14").Value = 'The Short Media is a marketing agency specializing in short-form video and social commerce. Founded by former TikTok leaders who shaped the platform''s advertising landscape, we bring a unique blend of innovation and expertise to our clients. Our team has a proven track record of driving exceptional results, having managed substantial ad budgets and developed groundbreaking strategies. We empower brands to connect with their audiences on TikTok, Meta, Snap, and social e-commerce, ensuring they stay ahead in the ever-evolving digital marketplace.'
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 'Hiring since March 2025'
$ws.Range("G14").Value = $null
$ws.Range("H14").Value = 'Data Analytics, Digital Marketing, Facebook Ads, Marketing, Job offer'
$ws.Range("I14").Value = 339
$ws.Range("J14").Value = 'Hyderabad'
$ws.Range("K14").Value = 'Stipend'
$ws.Range("L14").Value = '₹ 30,000 /month'
$ws.Range("M14").Value = $null
$ws.Range("N14").Value = 'https://internshala.com/internship/detail/performance-marketing-specialist-internship-in-hyderabad-at-the-short-media1741255003'

# Row 15: 'Medstown Private Limited'
$ws.Range("A15").Value = 'Medstown Private Limited'
$ws.Range("B15").Value = $null
$ws.Range("C15").Value = 'At Medstown our mission is to provide each customer with seamless and fast medicine delivery services, while not only saving their money and time but also uplifting and supporting our local pharmacies. We at Medstown aim to be the one stop shop and go to solution for patients in need of medicines, while ensuring our local pharmacies grow too. This way we aim to create a win-win situation for both customer and supplier.'
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 13
$ws.Range("F15").Value = 'Hiring since September 2023'
$ws.Range("G15").Value = $null
$ws.Range("H15").Value = 'Data Analysis, Effective Communication, Negotiation, Problem Solving, Time Management, Job offer'
$ws.Range("I15").Value = 177
$ws.Range("J15").Value = 'Hyderabad'
$ws.Range("K15").Value = 'Stipend'
$ws.Range("L15").Value = '₹ 10,000 /month'
$ws.Range("M15").Value = $null
$ws.Range("N15").Value = 'https://internshala.com/internship/detail/customer-success-internship-in-hyderabad-at-medstown-private-limited1741005176'

